$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text: every "Ready for handoff" cell becomes the handback message
#    (Overview!B2:B3 / C2:C3, zh-cn!C2:C3, de-de!C2:C3 all share this string).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2) Handback report: fill "Latest Target File" (F) and "Latest Handback
#    File" (G) with the same targets already linked from columns A/D, and
#    stamp the real "Latest Handback DateTime" (H) in place of the
#    placeholder "0001-01-01 00:00:00".
# ---------------------------------------------------------------------------

# --- zh-cn sheet ---
$zhMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/d076149e3e2aa5848457f2a2882088f694e5c8db/e2e/4112924c-117f-44a5-a62a-30ff92fa38dd.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/df9d4e5ff3b17e833907db3c7df61d7e484788c5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4112924c-117f-44a5-a62a-30ff92fa38dd.67d1fd250664effb642e5878015e34e7e382aae0.zh-cn.xlf"
$zhMdUrl2  = "https://github.com/OpenLocalizationTest/oltest/blob/d076149e3e2aa5848457f2a2882088f694e5c8db/e2e/ad4d94f7-7774-45f7-89f4-97867c76012d.md"
$zhXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/df9d4e5ff3b17e833907db3c7df61d7e484788c5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ad4d94f7-7774-45f7-89f4-97867c76012d.126fa422f86f4c1e0754461ae32b833be74e221f.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhMdUrl, "", "", "4112924c-117f-44a5-a62a-30ff92fa38dd.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfUrl, "", "", "4112924c-117f-44a5-a62a-30ff92fa38dd.67d1fd250664effb642e5878015e34e7e382aae0.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhMdUrl2, "", "", "ad4d94f7-7774-45f7-89f4-97867c76012d.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfUrl2, "", "", "ad4d94f7-7774-45f7-89f4-97867c76012d.126fa422f86f4c1e0754461ae32b833be74e221f.zh-cn.xlf")

$wsZh.Range("F2").Font.Underline = 2
$wsZh.Range("G2").Font.Underline = 2
$wsZh.Range("F3").Font.Underline = 2
$wsZh.Range("G3").Font.Underline = 2
$wsZh.Range("F2:G3").Font.Color = 15570276

# zh-cn: "Latest Handback DateTime" (H) gets a real timestamp now
$wsZh.Range("H2").Value = "2016-03-13 00:54:23"
$wsZh.Range("H3").Value = "2016-03-13 00:54:23"

# --- de-de sheet ---
$deMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/d076149e3e2aa5848457f2a2882088f694e5c8db/e2e/4112924c-117f-44a5-a62a-30ff92fa38dd.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e2cda626e6ff416ed4e23a1623d15a357303a0bd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4112924c-117f-44a5-a62a-30ff92fa38dd.67d1fd250664effb642e5878015e34e7e382aae0.de-de.xlf"
$deMdUrl2  = "https://github.com/OpenLocalizationTest/oltest/blob/d076149e3e2aa5848457f2a2882088f694e5c8db/e2e/ad4d94f7-7774-45f7-89f4-97867c76012d.md"
$deXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e2cda626e6ff416ed4e23a1623d15a357303a0bd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ad4d94f7-7774-45f7-89f4-97867c76012d.126fa422f86f4c1e0754461ae32b833be74e221f.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deMdUrl, "", "", "4112924c-117f-44a5-a62a-30ff92fa38dd.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfUrl, "", "", "4112924c-117f-44a5-a62a-30ff92fa38dd.67d1fd250664effb642e5878015e34e7e382aae0.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deMdUrl2, "", "", "ad4d94f7-7774-45f7-89f4-97867c76012d.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfUrl2, "", "", "ad4d94f7-7774-45f7-89f4-97867c76012d.126fa422f86f4c1e0754461ae32b833be74e221f.de-de.xlf")

$wsDe.Range("F2").Font.Underline = 2
$wsDe.Range("G2").Font.Underline = 2
$wsDe.Range("F3").Font.Underline = 2
$wsDe.Range("G3").Font.Underline = 2
$wsDe.Range("F2:G3").Font.Color = 15570276

# de-de: "Latest Handback DateTime" (H) gets its own real timestamp
$wsDe.Range("H2").Value = "2016-03-13 00:54:29"
$wsDe.Range("H3").Value = "2016-03-13 00:54:29"
